$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 23, pushing the existing rows 23-42 down to 24-43.
[void]$ws.Rows.Item(23).Insert()

# Populate the new row with the "PrimaryRelationship" extended-property
# definition (role-playing dimensions work-in-progress).
$ws.Range("A23").Value = "PrimaryRelationship"
$ws.Range("B23").Value = "Column"
$ws.Range("C23").Value = $true
$ws.Range("D23").Value = "Ordering of columns where multiple rol playing dimensions exist. E.g. Which table will also be used for DimUser as well as DimApprovingUser"

# Match the author's last selection in the saved file.
[void]$ws.Range("C24").Select()
